$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off (Y -> N) the Runmode flag for these suites
$ws.Range("C12").Value = "N"
$ws.Range("C16").Value = "N"
$ws.Range("C17").Value = "N"
$ws.Range("C20").Value = "N"
$ws.Range("C21").Value = "N"

# Update the active selection to C12
$ws.Range("C12").Select()
